$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 221.28572
$ws.Range("I9").Value = 129.8
$ws.Range("J9").Value = 450
$ws.Range("K9").Value = 129.8
$ws.Range("L9").Value = 450
$ws.Range("M9").Value = 39.19999999999999
$ws.Range("N9").Value = -788

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 72.5
$ws.Range("I33").Value = 72.5
$ws.Range("K33").Value = 72.5
$ws.Range("M33").Value = 156.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 646.125
$ws.Range("I38").Value = 234
$ws.Range("K38").Value = 702
$ws.Range("M38").Value = -330

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6228.8887
$ws.Range("I116").Value = 2927
$ws.Range("K116").Value = 2927
$ws.Range("M116").Value = 515

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1163.443
$ws.Range("J129").Value = 1171.9487
$ws.Range("L129").Value = 3515.8461
$ws.Range("N129").Value = -13515.8461

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3720.818
$ws.Range("I132").Value = 4109.1177
$ws.Range("J132").Value = 2400.6
$ws.Range("K132").Value = 12327.3531
$ws.Range("L132").Value = 7201.799999999999
$ws.Range("M132").Value = -9797.3531
$ws.Range("N132").Value = -12261.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1345.6444
$ws.Range("I138").Value = 547.2759
$ws.Range("J138").Value = 2792.6875
$ws.Range("K138").Value = 1641.8277
$ws.Range("L138").Value = 8378.0625
$ws.Range("M138").Value = 3498.1723
$ws.Range("N138").Value = -18658.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3647.0513
$ws.Range("I45").Value = 3312.389
$ws.Range("J45").Value = 3933.9048
$ws.Range("K45").Value = 3312.389
$ws.Range("L45").Value = 3933.9048
$ws.Range("M45").Value = -2935.389
$ws.Range("N45").Value = -4687.9048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2721.6
$ws.Range("I122").Value = 2721.6
$ws.Range("K122").Value = 8164.799999999999
$ws.Range("M122").Value = -5714.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 17002.484
$ws.Range("I132").Value = 1734.6428
$ws.Range("K132").Value = 5203.928400000001
$ws.Range("M132").Value = -2673.928400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5559645.5
$ws.Range("I105").Value = 4933.3335
$ws.Range("J105").Value = 16669070
$ws.Range("K105").Value = 4933.3335
$ws.Range("L105").Value = 16669070
$ws.Range("M105").Value = -3186.3335
$ws.Range("N105").Value = -16672564

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 888.3125
$ws.Range("I107").Value = 442.85715
$ws.Range("K107").Value = 442.85715
$ws.Range("M107").Value = 1477.14285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 32118.428
$ws.Range("I134").Value = 38480.38
$ws.Range("J134").Value = 1369
$ws.Range("K134").Value = 115441.14
$ws.Range("L134").Value = 4107
$ws.Range("M134").Value = -112906.14
$ws.Range("N134").Value = -9177

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 69917.5
$ws.Range("J68").Value = 69917.5
$ws.Range("L68").Value = 69917.5
$ws.Range("N68").Value = -71415.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 69917.5
$ws.Range("J71").Value = 69917.5
$ws.Range("L71").Value = 209752.5
$ws.Range("N71").Value = -217240.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1567.8
$ws.Range("I107").Value = 1069.625
$ws.Range("J107").Value = 2137.1428
$ws.Range("K107").Value = 1069.625
$ws.Range("L107").Value = 2137.1428
$ws.Range("M107").Value = 850.375
$ws.Range("N107").Value = -5977.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 35193.875
$ws.Range("I132").Value = 38292.07
$ws.Range("J132").Value = 13506.5
$ws.Range("K132").Value = 114876.21
$ws.Range("L132").Value = 40519.5
$ws.Range("M132").Value = -112346.21
$ws.Range("N132").Value = -45579.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 596.3333
$ws.Range("I122").Value = 316
$ws.Range("J122").Value = 1998
$ws.Range("K122").Value = 2844
$ws.Range("L122").Value = 17982
$ws.Range("M122").Value = -394
$ws.Range("N122").Value = -22882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 808.67
$ws.Range("I131").Value = 600
$ws.Range("J131").Value = 812.9286
$ws.Range("K131").Value = 1800
$ws.Range("L131").Value = 2438.7858
$ws.Range("M131").Value = 3240
$ws.Range("N131").Value = -12518.7858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2266
$ws.Range("I139").Value = 1502.5
$ws.Range("K139").Value = 4507.5
$ws.Range("M139").Value = 632.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4820037.5
$ws.Range("I70").Value = 24900
$ws.Range("J70").Value = 7816999
$ws.Range("K70").Value = 24900
$ws.Range("L70").Value = 7816999
$ws.Range("M70").Value = -24630
$ws.Range("N70").Value = -7817539

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4820037.5
$ws.Range("I73").Value = 24900
$ws.Range("J73").Value = 7816999
$ws.Range("K73").Value = 24900
$ws.Range("L73").Value = 7816999
$ws.Range("M73").Value = -23964
$ws.Range("N73").Value = -7818871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2679
$ws.Range("I97").Value = 1107.25
$ws.Range("J97").Value = 6870.3335
$ws.Range("K97").Value = 1107.25
$ws.Range("L97").Value = 6870.3335
$ws.Range("M97").Value = -611.25
$ws.Range("N97").Value = -7862.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2082.5806
$ws.Range("I113").Value = 1734
$ws.Range("J113").Value = 2302.7368
$ws.Range("K113").Value = 1734
$ws.Range("L113").Value = 2302.7368
$ws.Range("M113").Value = 436
$ws.Range("N113").Value = -6642.736800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2809.3333
$ws.Range("I122").Value = 2233.3333
$ws.Range("K122").Value = 6699.999899999999
$ws.Range("M122").Value = -4249.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 61413.54
$ws.Range("I132").Value = 51464.477
$ws.Range("K132").Value = 154393.431
$ws.Range("M132").Value = -151863.431

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3224.7273
$ws.Range("I7").Value = 3060.2104
$ws.Range("J7").Value = 4266.6665
$ws.Range("K7").Value = 3060.2104
$ws.Range("L7").Value = 4266.6665
$ws.Range("M7").Value = -2948.2104
$ws.Range("N7").Value = -4490.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 385.53333
$ws.Range("I16").Value = 377.35715
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 377.35715
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -207.35715
$ws.Range("N16").Value = -840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3224.7273
$ws.Range("I126").Value = 3060.2104
$ws.Range("J126").Value = 4266.6665
$ws.Range("K126").Value = 9180.6312
$ws.Range("L126").Value = 12799.9995
$ws.Range("M126").Value = -6710.6312
$ws.Range("N126").Value = -17739.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2110
$ws.Range("I132").Value = 1474
$ws.Range("K132").Value = 4422
$ws.Range("M132").Value = -1892

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2280.1
$ws.Range("I132").Value = 1114.8572
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 3344.5716
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -814.5715999999998
$ws.Range("N132").Value = -20057
